# Add a new weekly data row at the top of the data (row 2), pushing all
# existing data rows down by one (2..80 -> 3..81). This mirrors the source
# data feed being updated with a new week's record for "Caqui" (variety
# Mankaki, quality Primera) while keeping all previously recorded rows
# intact, just shifted down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2 (current data rows 2..80 -> 3..81)
$ws.Range("A2").EntireRow.Insert()

# The freshly inserted row inherits formatting from the row above (the
# bold header row); clear that so it looks like a normal data row, then
# reapply only the date number format used by the rest of column D.
$ws.Range("A2:T2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 2 keeps the same constant columns as the rest of the dataset
# (Mercado ID, Mercado, Región, Codreg, Tipo, Producto ID, Producto,
#  Categoría ID, Categoría, Calidad, Origen) and gets new values for the
# date, variety, volume, prices, unit and kg/unit columns.
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C2").Value = "Metropolitana"
$ws.Range("D2").Value = 45092
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = "Otros"
$ws.Range("I2").Value = 100107001
$ws.Range("J2").Value = "Caqui"
$ws.Range("K2").Value = "Mankaki"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 480
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 10500
$ws.Range("P2").Value = 10208
$ws.Range("Q2").Value = "`$/caja 15 kilos granel"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 681
$ws.Range("T2").Value = 15
